# Quiz data update: remove the "falls / 폭포" sound quiz entry, and replace the
# "mushroom / 버섯" entry with a new "bush / 덤불" entry (keeping its existing
# formatting), per commit "feat. sound/ effect ect".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the whole "falls" (폭포) quiz row (row 4). This shifts every
#    subsequent row up by one, so the old "mushroom" row (10) becomes row 9
#    and the old "owl" row (11) becomes row 10.
$ws.Rows("4:4").Delete()

# 2. The old "mushroom" entry (now row 9) is replaced in place with the new
#    "bush" entry, keeping the same problem text / sound path / formatting.
$ws.Range("B9").Value = "bush"
$ws.Range("D9").Value = "bush"
$ws.Range("F9").Value = "덤불"

# 3. Reflect the new selection recorded for the sheet.
$ws.Range("F10").Select()
